$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: Validation -> F1 train
$ws.Range("O1").Value = "F1 train"

# Row 2
$ws.Range("O2").Value = 0.810126582278481

# Row 3
$ws.Range("O3").Value = 1

# Row 4
$ws.Range("O4").Value = 1

# Row 5
$ws.Range("O5").Value = 0.8831168831168831

# Row 6 (MLP params + confusion matrix + metrics changed)
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 0.65
$ws.Range("J6").Value = 0.631578947368421
$ws.Range("K6").Value = 0.6
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.7
$ws.Range("N6").Value = 0.6
$ws.Range("O6").Value = 0.6052631578947368

# Row 7
$ws.Range("O7").Value = 1

# Row 8
$ws.Range("O8").Value = 1

# Row 9
$ws.Range("O9").Value = 1

# Row 10
$ws.Range("O10").Value = 0.918918918918919

# Row 11 (MLP params + confusion matrix + metrics changed)
$ws.Range("C11").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 5
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 0.55
$ws.Range("J11").Value = 0.5263157894736842
$ws.Range("K11").Value = 0.5
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("N11").Value = 0.5
$ws.Range("O11").Value = 0.9333333333333333

# Row 12
$ws.Range("O12").Value = 1

# Row 13
$ws.Range("O13").Value = 0.9866666666666667

# Row 14
$ws.Range("O14").Value = 1

# Row 15
$ws.Range("O15").Value = 0.7733333333333333

# Row 16
$ws.Range("O16").Value = 0.8051948051948052
